$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Col14a1", "Cd44", "ECs", 1, 0.3333333333333333, 0.2548253333333333, 0.764476, 0.00181063274871081, 0.00181063274871081, 3, 1, 261.380203, 784.1406089999999, 0.6968677182772199, 0.6968677182772199, 66.60629735620934, 599.4566762058839, 0.001261771512232113, 0.001261771512232113),
    @("ECs", "Col14a1", "Cd44", "FAPs", 1, 0.3333333333333333, 0.2548253333333333, 0.764476, 0.00181063274871081, 0.00181063274871081, 3, 1, 31.999428, 95.998284, 0.08531391482826334, 0.08531391482826335, 8.154264906576, 73.388384159184, 0.0001544721681087784, 0.0001544721681087784),
    @("ECs", "Col14a1", "Cd44", "sCs", 1, 0.3333333333333333, 0.2548253333333333, 0.764476, 0.00181063274871081, 0.00181063274871081, 3, 1, 81.69901900000001, 245.097057, 0.2178183668945166, 0.2178183668945167, 20.81897974968134, 187.370817747132, 0.0003943890683699184, 0.0003943890683699184),
    @("FAPs", "Col14a1", "Cd44", "ECs", 3, 1, 139.2710926666666, 417.813278, 0.9895750867169325, 0.9895750867169325, 3, 1, 261.380203, 784.1406089999999, 0.6968677182772199, 0.6968677182772199, 36402.70647324514, 327624.3582592062, 0.6896029327444108, 0.6896029327444108),
    @("FAPs", "Col14a1", "Cd44", "FAPs", 3, 1, 139.2710926666666, 417.813278, 0.9895750867169325, 0.9895750867169325, 3, 1, 31.999428, 95.998284, 0.08531391482826334, 0.08531391482826335, 4456.595302268327, 40109.35772041495, 0.08442452466433968, 0.0844245246643397),
    @("FAPs", "Col14a1", "Cd44", "sCs", 3, 1, 139.2710926666666, 417.813278, 0.9895750867169325, 0.9895750867169325, 3, 1, 81.69901900000001, 245.097057, 0.2178183668945166, 0.2178183668945167, 11378.31164592476, 102404.8048133228, 0.2155476293081819, 0.2155476293081819),
    @("sCs", "Col14a1", "Cd44", "ECs", 3, 1, 1.212359, 3.637077, 0.008614280534356691, 0.008614280534356693, 3, 1, 261.380203, 784.1406089999999, 0.6968677182772199, 0.6968677182772199, 316.886641528877, 2851.979773759893, 0.006003014020577018, 0.006003014020577019),
    @("sCs", "Col14a1", "Cd44", "FAPs", 3, 1, 1.212359, 3.637077, 0.008614280534356691, 0.008614280534356693, 3, 1, 31.999428, 95.998284, 0.08531391482826334, 0.08531391482826335, 38.794794530652, 349.153150775868, 0.0007349179958148736, 0.0007349179958148738),
    @("sCs", "Col14a1", "Cd44", "sCs", 3, 1, 1.212359, 3.637077, 0.008614280534356691, 0.008614280534356693, 3, 1, 81.69901900000001, 245.097057, 0.2178183668945166, 0.2178183668945167, 99.04854097582101, 891.4368687823891, 0.001876348517964799, 0.001876348517964799)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
